$wb = $excel.ActiveWorkbook

# A sheet that already has the "header row + column-A index" styling (style index 2)
# we want to reuse for the new/renamed sheets below.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# Step 1: the existing "总计" sheet becomes the new "2022-Q1" quarter sheet
# (keeps its original sheetId / position-derived rId "slot"); it is cleared and
# rebuilt with the fund-holding-detail layout used by the other quarter sheets.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1DataRng = $q1.Range("B2:G2")
$q1DataRng.NumberFormat = "@"
$q1.Range("B2").Value = "513030"
$q1.Range("C2").Value = "华安国际龙头(DAX)ETFQDII"
$q1.Range("D2").Value = "6.49"
$q1.Range("E2").Value = "92.80"
$q1.Range("F2").Value = "4.27"
$q1.Range("G2").Value = "0.2771"
$q1DataRng.ClearFormats()

$q1.Range("H2").Value = 7

$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: a brand-new "总计" sheet is created right after "2022-Q1" (gets a
# fresh sheetId) holding the refreshed summary table (quarters shifted down,
# 2022-Q1 row added on top).
# ---------------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"
$total.Cells.Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.28

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.22

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.23

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.36

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.41

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 0.46

$styleSrc.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
